$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values, keyed by row number, per the regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")
$gValues = @{
    2 = 4
    3 = 4
    4 = 1
    5 = 2
    6 = 1
    7 = 0
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 3
    22 = 1
    23 = 2
    24 = 3
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 2
    31 = 3
    32 = 0
    33 = 2
    34 = 2
    35 = 2
    36 = 2
    37 = 2
    38 = 2
    39 = 2
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 1
    45 = 2
    46 = 3
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 2
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 0
    57 = 0
    58 = 1
    59 = 2
    60 = 0
    61 = 1
    62 = 2
    63 = 2
    64 = 1
    65 = 3
    66 = 0
    67 = 3
    68 = 0
    69 = 0
    70 = 1
    71 = 0
    72 = 2
    73 = 1
    74 = 0
    75 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}

Write-Output "Updated $($gValues.Keys.Count) cells in column G"
